# Update cryptocurrency price (D) and 1h volume change (E) values
# Uses a leading apostrophe to force text entry (prevents Excel from
# auto-converting numeric-looking strings to floating point values),
# then resets the cell style so no stray quote-prefix style is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $range = $ws.Range($cellRef)
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

Set-TextValue "D2" "44.598.31"
Set-TextValue "E2" "  +0.64%  "
Set-TextValue "D3" "2.242.62"
Set-TextValue "E3" "  -0.30%  "
Set-TextValue "E4" "  +0.29%  "
Set-TextValue "D5" "305.87"
Set-TextValue "E5" "  -0.52%  "
Set-TextValue "D6" "95.12"
Set-TextValue "E6" "  -1.68%  "
Set-TextValue "E7" "  -0.58%  "
Set-TextValue "E8" "  +0.12%  "
Set-TextValue "E9" "  -1.75%  "
Set-TextValue "D10" "35.47"
Set-TextValue "E10" "  +0.50%  "
Set-TextValue "E11" "  -1.61%  "
Set-TextValue "D12" "7.23"
Set-TextValue "E12" "  -0.93%  "
Set-TextValue "E13" "  -0.24%  "
Set-TextValue "D14" "2.583.96"
Set-TextValue "E14" "  -0.36%  "
Set-TextValue "D15" "2.231.11"
Set-TextValue "E15" "  -4.40%  "
Set-TextValue "D16" "0.832"
Set-TextValue "E16" "  -0.61%  "
Set-TextValue "D17" "13.57"
Set-TextValue "E17" "  -0.69%  "
Set-TextValue "D18" "44.391.49"
Set-TextValue "E18" "  +0.60%  "
Set-TextValue "D19" "0.0₃0941"
Set-TextValue "E19" "  -3.16%  "
Set-TextValue "D20" "11.86"
Set-TextValue "E20" "  -2.71%  "
Set-TextValue "D21" "6.19"
Set-TextValue "E21" "  -3.23%  "
Set-TextValue "E22" "  -0.75%  "
Set-TextValue "D23" "236.94"
Set-TextValue "E23" "  -0.52%  "
Set-TextValue "E24" "  -0.85%  "
Set-TextValue "E25" "  -1.65%  "
Set-TextValue "E26" "  -0.14%  "
Set-TextValue "E27" "  +8.01%  "
Set-TextValue "D28" "9.76"
Set-TextValue "D29" "37.17"
Set-TextValue "E29" "  -4.58%  "
Set-TextValue "E30" "  -0.65%  "
Set-TextValue "D31" "19.90"
Set-TextValue "E31" "  -1.15%  "
Set-TextValue "D32" "149.78"
Set-TextValue "E32" "  -1.72%  "
Set-TextValue "D33" "0.0786"
Set-TextValue "E33" "  -1.48%  "
Set-TextValue "E34" "  -0.17%  "
Set-TextValue "E35" "  -3.13%  "
Set-TextValue "E36" "  +0.57%  "
Set-TextValue "E37" "  -1.55%  "
Set-TextValue "E38" "  +5.59%  "
Set-TextValue "D39" "15.22"
Set-TextValue "E39" "  +3.64%  "
Set-TextValue "D40" "3.40"
Set-TextValue "E40" "  -6.29%  "
Set-TextValue "D41" "3.78"
Set-TextValue "E41" "  -2.22%  "
Set-TextValue "D42" "0.0299"
Set-TextValue "E42" "  -0.43%  "
Set-TextValue "E43" "  +0.09%  "
Set-TextValue "D44" "1.813.12"
Set-TextValue "E44" "  +3.20%  "
Set-TextValue "E45" "  +11.62%  "
Set-TextValue "D46" "81.47"
Set-TextValue "E46" "  -2.29%  "
Set-TextValue "D47" "0.189"
Set-TextValue "E47" "  -1.89%  "
Set-TextValue "D48" "98.38"
Set-TextValue "E48" "  -2.25%  "
Set-TextValue "E49" "  -2.74%  "
Set-TextValue "D50" "68.51"
Set-TextValue "E50" "  +0.65%  "
Set-TextValue "D51" "54.27"
Set-TextValue "E51" "  -1.33%  "
